$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing row (141) into each new row
$ws.Range("A141:AC141").Copy()
$ws.Range("A142:AC145").PasteSpecial(-4122)

# Row 142
$ws.Range("A142").Value = 140
$ws.Range("B142").Value = 7011634
$ws.Range("C142").Value = "Azerbaijan Premier League"
$ws.Range("D142").Value = "Azerbaijan Premier League"
$ws.Range("E142").Value = 45388.375
$ws.Range("F142").Value = "PFK Turan Tovuz"
$ws.Range("G142").Value = "Sabah"
$ws.Range("H142").Value = 2
$ws.Range("I142").Value = 0
$ws.Range("J142").Value = "H"
$ws.Range("K142").Value = 2.55
$ws.Range("L142").Value = 3.25
$ws.Range("M142").Value = 2.4
$ws.Range("N142").Value = 2.625
$ws.Range("O142").Value = 3.2
$ws.Range("P142").Value = 2.375
$ws.Range("Q142").Value = 0
$ws.Range("R142").Value = 2
$ws.Range("S142").Value = 1.8
$ws.Range("T142").Value = 2.25
$ws.Range("U142").Value = 1.8
$ws.Range("V142").Value = 2
$ws.Range("W142").Value = 1.625
$ws.Range("X142").Value = -1
$ws.Range("Y142").Value = -1
$ws.Range("Z142").Value = 1
$ws.Range("AA142").Value = -1
$ws.Range("AB142").Value = -0.5
$ws.Range("AC142").Value = 0.5

# Row 143
$ws.Range("A143").Value = 141
$ws.Range("B143").Value = 7011633
$ws.Range("C143").Value = "Azerbaijan Premier League"
$ws.Range("D143").Value = "Azerbaijan Premier League"
$ws.Range("E143").Value = 45388.54166666666
$ws.Range("F143").Value = "FK Qarabag"
$ws.Range("G143").Value = "Sabail FC"
$ws.Range("H143").Value = 4
$ws.Range("I143").Value = 2
$ws.Range("J143").Value = "H"
$ws.Range("K143").Value = 1.222
$ws.Range("L143").Value = 5.5
$ws.Range("M143").Value = 8.5
$ws.Range("N143").Value = 1.181
$ws.Range("O143").Value = 7
$ws.Range("P143").Value = 9
$ws.Range("Q143").Value = -2.25
$ws.Range("R143").Value = 1.95
$ws.Range("S143").Value = 1.85
$ws.Range("T143").Value = 3.75
$ws.Range("U143").Value = 1.975
$ws.Range("V143").Value = 1.825
$ws.Range("W143").Value = 0.181
$ws.Range("X143").Value = -1
$ws.Range("Y143").Value = -1
$ws.Range("Z143").Value = -0.5
$ws.Range("AA143").Value = 0.425
$ws.Range("AB143").Value = 0.9750000000000001
$ws.Range("AC143").Value = -1

# Row 144
$ws.Range("A144").Value = 142
$ws.Range("B144").Value = 7011635
$ws.Range("C144").Value = "Azerbaijan Premier League"
$ws.Range("D144").Value = "Azerbaijan Premier League"
$ws.Range("E144").Value = 45389.39583333334
$ws.Range("F144").Value = "Zira IK"
$ws.Range("G144").Value = "FK Gabala"
$ws.Range("H144").Value = 4
$ws.Range("I144").Value = 0
$ws.Range("J144").Value = "H"
$ws.Range("K144").Value = 1.4
$ws.Range("L144").Value = 4.333
$ws.Range("M144").Value = 6
$ws.Range("N144").Value = 1.727
$ws.Range("O144").Value = 3.75
$ws.Range("P144").Value = 3.75
$ws.Range("Q144").Value = -0.5
$ws.Range("R144").Value = 1.75
$ws.Range("S144").Value = 1.95
$ws.Range("T144").Value = 2.25
$ws.Range("U144").Value = 1.9
$ws.Range("V144").Value = 1.9
$ws.Range("W144").Value = 0.7270000000000001
$ws.Range("X144").Value = -1
$ws.Range("Y144").Value = -1
$ws.Range("Z144").Value = 0.75
$ws.Range("AA144").Value = -1
$ws.Range("AB144").Value = 0.8999999999999999
$ws.Range("AC144").Value = -1

# Row 145
$ws.Range("A145").Value = 143
$ws.Range("B145").Value = 7011636
$ws.Range("C145").Value = "Azerbaijan Premier League"
$ws.Range("D145").Value = "Azerbaijan Premier League"
$ws.Range("E145").Value = 45389.54166666666
$ws.Range("F145").Value = "Neftchi Baku"
$ws.Range("G145").Value = "FK Sumqayit"
$ws.Range("H145").Value = 1
$ws.Range("I145").Value = 1
$ws.Range("J145").Value = "D"
$ws.Range("K145").Value = 1.727
$ws.Range("L145").Value = 3.6
$ws.Range("M145").Value = 4
$ws.Range("N145").Value = 1.727
$ws.Range("O145").Value = 3.6
$ws.Range("P145").Value = 4
$ws.Range("Q145").Value = -0.5
$ws.Range("R145").Value = 1.75
$ws.Range("S145").Value = 1.95
$ws.Range("T145").Value = 2.25
$ws.Range("U145").Value = 1.8
$ws.Range("V145").Value = 2
$ws.Range("W145").Value = -1
$ws.Range("X145").Value = 2.6
$ws.Range("Y145").Value = -1
$ws.Range("Z145").Value = -1
$ws.Range("AA145").Value = 0.95
$ws.Range("AB145").Value = -0.5
$ws.Range("AC145").Value = 0.5

